# Adds the "vCPU" and "RAM" horizontal/vertical scaling threshold sheets
# (UpperThresholdValue, VcpuValueForHScaling, VcpuvalueForVScaling,
#  ToCheckAlert, RAMValueForHScaling, RAMValueForVScaling) after the
# existing LowerThresholdValue sheet, matching the "Till RAM of both
# Horizontal And Vertical" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: build a simple two-/three-column "Sl.no" lookup sheet of the
# shape every sheet in this workbook already uses:
#   row1: headers, row2: Sl.no=1 plus one or two text-formatted values
# (NOTE: positional params only -- named "-Param value" binding is not
#  reliable against this host, so every call below passes args in order)
# ---------------------------------------------------------------------
function New-LookupSheet {
    param(
        [string]$Name,
        [string]$HeaderA,
        [string]$HeaderB,
        [string]$ValueB,
        [string]$HeaderC,
        [string]$ValueC
    )

    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $after)
    $ws.Name = $Name

    $ws.Range("A1").Value = $HeaderA
    $ws.Range("B1").Value = $HeaderB
    $ws.Range("A2").Value = 1

    $ws.Range("B2").NumberFormat = "@"
    $ws.Range("B2").Value = $ValueB

    if ($HeaderC) {
        $ws.Range("C1").Value = $HeaderC
        $ws.Range("C2").NumberFormat = "@"
        $ws.Range("C2").Value = $ValueC
    }

    return $ws
}

# --- sheet7: UpperThresholdValue --------------------------------------
$wsUpper = New-LookupSheet "UpperThresholdValue" "Sl.no" "UpperThresholdValue" "80" $null $null
$wsUpper.Range("A1:B2").Select()

# --- sheet8: VcpuValueForHScaling -------------------------------------
$wsVcpuH = New-LookupSheet "VcpuValueForHScaling" "Sl.no" "VcpuValueForHScaling" "6" $null $null
$wsVcpuH.Range("A1:B2").Select()

# --- sheet9: VcpuvalueForVScaling -------------------------------------
$wsVcpuV = New-LookupSheet "VcpuvalueForVScaling" "Sl.no" "MaxVcpuValueForVScaling" "16" "MinVcpuValueForVScaling" "14"
$wsVcpuV.Range("A1:C2").Select()

# --- sheet10: ToCheckAlert ---------------------------------------------
$wsAlert = New-LookupSheet "ToCheckAlert" "Sl.no(min value should be greater than max value)" "MaxVcpuValueForVScaling" "14" "MinVcpuValueForVScaling" "16"
$wsAlert.Range("C1").Select()

# --- sheet11: RAMValueForHScaling --------------------------------------
$wsRamH = New-LookupSheet "RAMValueForHScaling" "Sl.no" "RamValueForHScaling" "100" $null $null
$wsRamH.Range("B2").Select()

# --- sheet12: RAMValueForVScaling --------------------------------------
$wsRamV = New-LookupSheet "RAMValueForVScaling" "Sl.no" "MaxRAMValueForVScaling" "2" "MinRAMValueForVScaling" "14"
$wsRamV.Range("B2").Select()

# ---------------------------------------------------------------------
# View-state touch-ups on the pre-existing sheets
# ---------------------------------------------------------------------

# Location: selection moves to B2 (was the whole A1:B2 block)
$wsLocation = $wb.Worksheets.Item("Location")
$wsLocation.Range("B2").Select()

# LowerThresholdValue: no longer the active tab; selection becomes the
# full data block instead of the single cell it used to be.
$wsLower = $wb.Worksheets.Item("LowerThresholdValue")
$wsLower.Range("A1:B2").Select()

# The newly added RAMValueForVScaling sheet ends up being the active /
# selected tab, matching the workbook's new activeTab.
$wsRamV.Activate()
$wsRamV.Range("B2").Select()
